{"js": "const body = context.document.body;\n\n// Update the date heading paragraph (first paragraph in the body)\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2023-11-11 Saturday\", \"Replace\");\n\n// Update each answer cell in the table, in row-major order (0-indexed)\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst values = [\n  \"37+19=56\",\n  \"98-2=96\",\n  \"31-26=5\",\n  \"93-77=16\",\n  \"64+18=82\",\n  \"11+81=92\",\n  \"11+55=66\",\n  \"44-22=22\",\n  \"0+88=88\",\n  \"11+52=63\",\n  \"35+30=65\",\n  \"23+51=74\",\n  \"81+9=90\",\n  \"35+10=45\",\n  \"18+34=52\",\n  \"65-52=13\",\n  \"86-66=20\",\n  \"63-32=31\",\n  \"67+10=77\",\n  \"67-59=8\",\n  \"60+32=92\",\n  \"63-16=47\",\n  \"64-25=39\",\n  \"94-63=31\",\n  \"78-43=35\",\n  \"83-6=77\",\n  \"64+2=66\",\n  \"11+58=69\",\n  \"22+75=97\",\n  \"59+39=98\",\n  \"88-57=31\",\n  \"30+52=82\",\n  \"3+19=22\",\n  \"52-43=9\",\n  \"55+3=58\",\n  \"13+51=64\",\n  \"60+38=98\",\n  \"48-42=6\",\n  \"22+36=58\",\n  \"20+12=32\",\n  \"95-19=76\",\n  \"92+4=96\",\n  \"33-12=21\",\n  \"82-64=18\",\n  \"78-64=14\",\n  \"42+43=85\",\n  \"88-11=77\",\n  \"18-7=11\",\n  \"2+32=34\",\n  \"0+75=75\",\n  \"44+33=77\",\n  \"80-10=70\",\n  \"90-40=50\",\n  \"65-26=39\",\n  \"90-31=59\",\n  \"36-6=30\",\n  \"74+19=93\",\n  \"61+28=89\",\n  \"64+3=67\",\n  \"80+4=84\",\n  \"66-23=43\",\n  \"3+10=13\",\n  \"97-52=45\",\n  \"84-61=23\",\n  \"54-32=22\",\n  \"26-17=9\",\n  \"66-25=41\",\n  \"77+9=86\",\n  \"43+20=63\",\n  \"39+23=62\",\n  \"26-15=11\",\n  \"84-43=41\",\n  \"42-2=40\",\n  \"7+16=23\",\n  \"6+69=75\",\n  \"59+10=69\",\n  \"0+20=20\",\n  \"36-16=20\",\n  \"89-15=74\",\n  \"8+73=81\",\n  \"99-14=85\",\n  \"69-47=22\",\n  \"8+47=55\",\n  \"41+38=79\",\n  \"99-58=41\",\n  \"88-31=57\",\n  \"38+30=68\",\n  \"60+26=86\",\n  \"65-28=37\",\n  \"31+3=34\",\n  \"47+1=48\",\n  \"50-46=4\",\n  \"28+16=44\",\n  \"73-18=55\",\n  \"40-7=33\",\n  \"12-3=9\",\n  \"76-4=72\",\n  \"19+69=88\",\n  \"43-31=12\",\n  \"88-72=16\"\n];\n\nconst columns = 5;\nfor (let i = 0; i < values.length; i++) {\n  const r = Math.floor(i / columns);\n  const c = i % columns;\n  table.getCell(r, c).value = values[i];\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading paragraph\n$p = $d.Paragraphs.Item(1)\n$pr = $p.Range\n$pr.End = $pr.End - 1\n$pr.Text = \"2023-11-11 Saturday\"\n\n# Update each answer cell in the table, in row-major order\n$values = @(\n  \"37+19=56\",\n  \"98-2=96\",\n  \"31-26=5\",\n  \"93-77=16\",\n  \"64+18=82\",\n  \"11+81=92\",\n  \"11+55=66\",\n  \"44-22=22\",\n  \"0+88=88\",\n  \"11+52=63\",\n  \"35+30=65\",\n  \"23+51=74\",\n  \"81+9=90\",\n  \"35+10=45\",\n  \"18+34=52\",\n  \"65-52=13\",\n  \"86-66=20\",\n  \"63-32=31\",\n  \"67+10=77\",\n  \"67-59=8\",\n  \"60+32=92\",\n  \"63-16=47\",\n  \"64-25=39\",\n  \"94-63=31\",\n  \"78-43=35\",\n  \"83-6=77\",\n  \"64+2=66\",\n  \"11+58=69\",\n  \"22+75=97\",\n  \"59+39=98\",\n  \"88-57=31\",\n  \"30+52=82\",\n  \"3+19=22\",\n  \"52-43=9\",\n  \"55+3=58\",\n  \"13+51=64\",\n  \"60+38=98\",\n  \"48-42=6\",\n  \"22+36=58\",\n  \"20+12=32\",\n  \"95-19=76\",\n  \"92+4=96\",\n  \"33-12=21\",\n  \"82-64=18\",\n  \"78-64=14\",\n  \"42+43=85\",\n  \"88-11=77\",\n  \"18-7=11\",\n  \"2+32=34\",\n  \"0+75=75\",\n  \"44+33=77\",\n  \"80-10=70\",\n  \"90-40=50\",\n  \"65-26=39\",\n  \"90-31=59\",\n  \"36-6=30\",\n  \"74+19=93\",\n  \"61+28=89\",\n  \"64+3=67\",\n  \"80+4=84\",\n  \"66-23=43\",\n  \"3+10=13\",\n  \"97-52=45\",\n  \"84-61=23\",\n  \"54-32=22\",\n  \"26-17=9\",\n  \"66-25=41\",\n  \"77+9=86\",\n  \"43+20=63\",\n  \"39+23=62\",\n  \"26-15=11\",\n  \"84-43=41\",\n  \"42-2=40\",\n  \"7+16=23\",\n  \"6+69=75\",\n  \"59+10=69\",\n  \"0+20=20\",\n  \"36-16=20\",\n  \"89-15=74\",\n  \"8+73=81\",\n  \"99-14=85\",\n  \"69-47=22\",\n  \"8+47=55\",\n  \"41+38=79\",\n  \"99-58=41\",\n  \"88-31=57\",\n  \"38+30=68\",\n  \"60+26=86\",\n  \"65-28=37\",\n  \"31+3=34\",\n  \"47+1=48\",\n  \"50-46=4\",\n  \"28+16=44\",\n  \"73-18=55\",\n  \"40-7=33\",\n  \"12-3=9\",\n  \"76-4=72\",\n  \"19+69=88\",\n  \"43-31=12\",\n  \"88-72=16\"\n)\n\n$t = $d.Tables.Item(1)\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cr = $cell.Range\n    $cr.End = $cr.End - 1\n    $cr.Text = $values[$idx]\n    $idx = $idx + 1\n  }\n}"}
